$wb = $excel.ActiveWorkbook
$template = $wb.Worksheets.Item("2021-Q4")

# New quarter sheet goes right after "2021-Q4" and before "总计"
$newSheet = $wb.Worksheets.Add($null, $template)
$newSheet.Name = "2022-Q1"

# Copy header style (bold/centered/bordered) from an existing fund-detail sheet
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Copy column-A numbering style too
$template.Range("A2").Copy()
$newSheet.Range("A2:A9").PasteSpecial(-4122)

# Force text format on the numeric-looking columns only, so values such as
# "013516" (leading zero) and "5.00" / "2.4400" (trailing zeros) are kept as
# literal text instead of being coerced to Double by the COM layer.
$newSheet.Range("B2:B9").NumberFormat = "@"
$newSheet.Range("D2:G9").NumberFormat = "@"

$data = @(
    @(0, "519066", "汇添富蓝筹稳健混合", "65.24", "74.02", "3.74", "2.4400", 8),
    @(1, "013516", "汇添富蓝筹稳健混合E", "65.24", "74.02", "3.74", "2.4400", 8),
    @(2, "257010", "国联安小盘精选混合", "9.15", "74.36", "4.42", "0.4044", 6),
    @(3, "550001", "信诚四季红混合", "5.00", "72.84", "4.12", "0.2060", 1),
    @(4, "002367", "国联安安稳灵活配置混合", "2.32", "33.99", "2.55", "0.0592", 4),
    @(5, "011284", "中信保诚龙腾精选混合", "1.22", "75.38", "4.22", "0.0515", 1),
    @(6, "006209", "中信保诚新蓝筹灵活配置混合", "1.16", "77.03", "4.21", "0.0488", 1),
    @(7, "006138", "国联安价值优选股票", "0.60", "93.30", "4.77", "0.0286", 5)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# --- Update the "总计" (summary) sheet: insert a new 2022-Q1 row at the top ---
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()
$summary.Rows.Item(2).ClearFormats()

$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q1"
$summary.Cells.Item(2, 3).Value = 8
$summary.Cells.Item(2, 4).Value = 5.68

# Renumber the row-index column (A) for the rows that got shifted down
for ($row = 3; $row -le 7; $row++) {
    $summary.Cells.Item($row, 1).Value = $row - 2
}

Write-Output "done"
